# The diary table cell currently reads (split across three runs,
# wrapped in spell-check proofErr markers around "Gantt"):
#   "Gantt e finire la doc in modo che posso iniziare a programmare "
# It needs to become a single, clean run reading:
#   "Fare il gantt e finire la parte di progettazione della doc"
# with the proofErr markers gone entirely.

$d = $word.ActiveDocument

$oldSnippet = "Gantt e finire la doc"
$newText    = "Fare il gantt e finire la parte di progettazione della doc"

# Locate the table cell that holds the paragraph we need to rewrite.
$targetCell = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    for ($rIdx = 1; $rIdx -le $tbl.Rows.Count; $rIdx++) {
        for ($cIdx = 1; $cIdx -le $tbl.Columns.Count; $cIdx++) {
            try {
                $c = $tbl.Cell($rIdx, $cIdx)
            } catch {
                continue
            }
            if ($c.Range.Text -like ("*" + $oldSnippet + "*")) {
                $targetCell = $c
            }
        }
    }
}

$p = $targetCell.Range.Paragraphs.Item(1)
$r = $p.Range

# Rebuild the paragraph from clean OOXML (single run, no proofErr marks)
# via InsertXML so the stray <w:proofErr/> bookmarks that wrapped the old
# "Gantt" text are dropped along with the old runs, instead of merely
# replacing the visible text and leaving orphaned proofErr markers behind.
$xmlFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3724A674" w14:textId="3300D6A0" w:rsidR="00632B06" w:rsidRDefault="00C50B4C" w:rsidP="00434F37"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Fare il gantt e finire la parte di progettazione della doc</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xmlFrag) | Out-Null
